$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a plain-numeric-looking string to be stored as TEXT
# (mirrors the source data's inlineStr cells) without leaving a lasting
# number-format override on the cell.
function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

$ws.Range("D2").Value = '26.220.28'
$ws.Range("E2").Value = '  +0.24%  '
$ws.Range("D3").Value = '1.604.22'
$ws.Range("E3").Value = '  +0.14%  '
$ws.Range("E4").Value = '  -0.03%  '
Set-TextValue "D5" '212.60'
$ws.Range("E5").Value = '  -0.16%  '
$ws.Range("E6").Value = '  -0.03%  '
$ws.Range("E7").Value = '  -0.23%  '
$ws.Range("E8").Value = '  +0.15%  '
$ws.Range("E9").Value = '  -0.49%  '
Set-TextValue "D10" '18.43'
$ws.Range("E10").Value = '  +2.43%  '
Set-TextValue "D11" '0.0815'
$ws.Range("E11").Value = '  -0.22%  '
$ws.Range("D12").Value = '1.827.38'
$ws.Range("E12").Value = '  +0.13%  '
$ws.Range("D13").Value = '1.612.42'
$ws.Range("E13").Value = '  +0.60%  '
$ws.Range("E14").Value = '  +0.41%  '
Set-TextValue "D15" '0.513'
$ws.Range("E15").Value = '  +0.30%  '
$ws.Range("D16").Value = '26.196.67'
$ws.Range("E16").Value = '  +0.20%  '
Set-TextValue "D17" '61.92'
$ws.Range("E17").Value = '  +2.28%  '
$ws.Range("E18").Value = '  +0.97%  '
$ws.Range("E19").Value = '  +0.01%  '
Set-TextValue "D20" '200.89'
$ws.Range("E20").Value = '  -2.17%  '
$ws.Range("E21").Value = '  +0.82%  '
Set-TextValue "D22" '9.30'
$ws.Range("E22").Value = '  -0.01%  '
$ws.Range("E23").Value = '  +0.05%  '
Set-TextValue "D24" '1.87'
$ws.Range("E24").Value = '  +2.25%  '
Set-TextValue "D25" '143.82'
$ws.Range("E25").Value = '  +1.18%  '
$ws.Range("E26").Value = '  +0.00%  '
$ws.Range("E27").Value = '  -1.93%  '
Set-TextValue "D28" '15.21'
$ws.Range("E28").Value = '  -0.05%  '
$ws.Range("E29").Value = '  +1.85%  '
$ws.Range("E30").Value = '  +3.96%  '
$ws.Range("E31").Value = '  +0.43%  '
$ws.Range("E32").Value = '  +2.33%  '
$ws.Range("E33").Value = '  -1.11%  '
$ws.Range("E34").Value = '  +0.51%  '
$ws.Range("E35").Value = '  +1.80%  '
$ws.Range("D36").Value = '1.162.16'
$ws.Range("E36").Value = '  +4.41%  '
$ws.Range("E37").Value = '  +3.35%  '
$ws.Range("E38").Value = '  -0.06%  '
$ws.Range("E39").Value = '  -0.34%  '
$ws.Range("E40").Value = '  +0.83%  '
$ws.Range("E41").Value = '  +0.73%  '
Set-TextValue "D42" '0.783'
$ws.Range("E42").Value = '  +0.33%  '
$ws.Range("E43").Value = '  +3.87%  '
$ws.Range("D44").Value = '1.738.82'
Set-TextValue "D45" '91.64'
$ws.Range("E45").Value = '  -1.28%  '
$ws.Range("D46").Value = '0.0₆0106'
$ws.Range("E46").Value = '  +19.19%  '
$ws.Range("E47").Value = '  +1.18%  '
$ws.Range("E48").Value = '  +1.33%  '
$ws.Range("E49").Value = '  +0.09%  '
$ws.Range("E50").Value = '  -0.57%  '
$ws.Range("E51").Value = '  -0.13%  '
